$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mark Yean Tuck Ming"
$ws.Range("C2").Value = "+6012-225 1051"
$ws.Range("D2").Value = "yeanmark@gmail.com"

$ws.Range("G2").Value = "[{'job_title': 'Data Engineer', 'job_company': 'CelcomDigi Sdn Bhd', 'Industries': 'Information and communication', 'start_date': '2022-05', 'end_date': '2023-05', 'job_location': 'N/A', 'job_duration': '1.0'}, {'job_title': 'Data Analyst Intern', 'job_company': 'MC Crenergy Sdn Bhd', 'Industries': 'Information and communication', 'start_date': '2019-11', 'end_date': '2020-03', 'job_location': 'N/A', 'job_duration': '0.3'}]"

$ws.Range("H2").Value = "[{'Country': 'Malaysia', 'State': 'Selangor', 'City': 'Puchong'}]"

$ws.Range("I2").Value = "[{'field_of_study': 'Intelligence System', 'level': ""Bachelor's Degree"", 'cgpa': 'N/A', 'university': 'Asia Pacific University', 'start_date': '2018-02', 'year_of_graduation': '2021'}, {'field_of_study': 'Information Technology', 'level': 'Foundation', 'cgpa': 'N/A', 'university': 'Asia Pacific University', 'start_date': '2017-02', 'year_of_graduation': '2017'}]"

$ws.Range("J2").Value = "['Project Management Workshop Series', 'Body Language Secrets: Nonverbal Communication Strategies', 'Computer Hacking Forensic Investigator', 'Certified Ethical Hacking Ver. 10', 'APIIT Certified Security Professional', 'Premier-Pride Challenge 2019', 'Artificial Intelligence Modules', 'AI and Machine Learning Competence for Industry 4.0', 'Certified Engineer in Computer Vision', 'Dataiku Core Designer Certificate']"

$ws.Range("K2").Value = "['Python', 'Selenium', 'BeautifulSoup', 'Pandas', 'Numpy', 'Matplotlib', 'Seaborn', 'Plotly', 'Supervised ML models', 'Unsupervised ML models', 'Natural Language Process', 'Computer Vision', 'Flask', 'DJango', 'Dash', 'MySQL', 'Teradata SQL', 'Oracle Data Integrator', 'Linux', 'C++', 'C#', 'Java', 'Microsoft Power BI', 'OLAP Cube', 'SSIS', 'ASP .net', 'DevExpress']"

$ws.Range("L2").Value = "['Chinese(Mandarin)', 'Chinese(Cantonese)', 'English', 'Malay']"
